$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update price in cell C8 from 15599 to 16499
$ws.Range("C8").Value = 16499

# Select cell C8 (reflects the updated selection in the sheet view)
$ws.Range("C8").Select()
